# Convert text numbers/bools to native types: add a new "isEmployee" row
# above the "phones[0].type" row, with formulas that render boolean-ish
# text ("true"/"false") as native string results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 ("phones[0].type"), shifting the existing
# phones/aliases rows down by one.
$ws.Rows(7).Insert() | Out-Null

# Populate the new row 7 with the isEmployee flag fields.
$ws.Range("A7").Value = "isEmployee"
$ws.Range("B7").Formula = '="true"'
$ws.Range("C7").Formula = '="false"'

# Match the saved selection state: row 8 ("phones[0].type") selected.
$ws.Rows(8).Select() | Out-Null
